$d = $word.ActiveDocument

# 1. Title (Heading1)
$d.Content.Find.Execute("Review 187: Transformers are RNNs: Fast Autoregressive Transformers with Linear Attention", $true, $false, $false, $false, $false, $true, 1, $false, "Review 186: HiPPO: Recurrent Memory with Optimal Polynomial Projections", 2) | Out-Null

# 2. Paper link (bold run)
$d.Content.Find.Execute("Paper: https://arxiv.org/abs/2006.16236v3", $true, $false, $false, $false, $false, $true, 1, $false, "Paper: https://arxiv.org/abs/2008.07669v2", 2) | Out-Null

# 3. pdf/abs link
$d.Content.Find.Execute("https://arxiv.org/pdf/2006.16236", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2008.07669", 2) | Out-Null

# 4. Paragraph (intro)
$d.Content.Find.Execute("אחרי הסקירה הקודמת הכבדה מאוד מחכה לנו היום סקירה קלילה (הסקירה הבאה הולכת להיות די כבדה). כמו שכבר אמרנו אחד החסרונות הבולטים של הטרנספורמר היא הסיבוכיות הריבועית שלו במונחי אורך הקלט (= מספר איברים בסדרת הקלט). הסיבוכיות הזו בא על ידי ביטוי גם במהלך האימון וגם במהלך ההיסק (inference). סיבוכיות ריבועית זאת כואבת במיוחד בזמן ההיסק כאשר אין לנו יכולת לחזות מספר טוקנים בו זמנית כי לחיזוי טוקן n אנו צריכים לדעת את ה-(n-1) הטוקנים הראשונים. האם ניתן להפוך את הטרנספורמר לסוג של RNN במהלך ההיסק כאשר כל הזיכרון על הטוקנים הקודמים נדחס לכמה וקטורים בודדים (וקטור זכרון ווקטור של המצב)?", $true, $false, $false, $false, $false, $true, 1, $false, "הגענו למאמר השני בסדרה - המאמר הזה חשוב מאוד כי הוא מפתח בסיס מתמטי מוצק המשמש כל המודלים מבוססים על מערכות דינמיות לינאריות כולל כמובן ממבה. המאמר הזה קצת (די הרבה) כבד מתמטית אך אנסה לעשות כמיטב יכולתי כדי להעביר לכם את המסר העיקרי שהוא מביא איתו. ", 2) | Out-Null

# 5. Paragraph (body 1)
$d.Content.Find.Execute("הטרנספורמר המקורי אינו מאפשר אופן חישוב כזה כי הוא מכיל פעולה לא לינארית (softmax) בתוך מנגנון תשומת הלב שלו. ניתן לראות די בקלות שלא ניתן לעקוף את מגבלת הסיבוכיות הריבועית שלו ללא שינוי של אופן חישוב של תשומת הלב. המאמר המסוקר מציע להחליף את חישוב הסופטמקס במנגנון זה בחישוב לינארי (מכפלת מטריצות) המחושבות על ידי הפעלת פונקציה לא לינארית phi על וקטורי השאילתות Q ושל וקטורי המפתחות K. מי שעוד זוכר מה זה KT)Kernel Trick) מבין מה שנעשה כאן הוא KT בכיוון ההפוך. ", $true, $false, $false, $false, $false, $true, 1, $false, "בסקירה הקודמת דיברנו על איך ניתן לבנות וקטור זיכרון (m(t בעל יכולת לשחזר פונקצית קלט (u(x ל-x מאינטרוול ; כאן t מסמן גודל חלון הקשר (כלומר אורך הזיכרון). פונקצית (m(t ממודלת על ידי מערכת דינמית לינארית ושילובה עם פולינומי Legendre משחזר לנו את הקלט u. נעיר שאנו עובדים עם הגרסאות הדיסקרטיות של המודלים האלו שהן בעצם נוסחת נסיגה עבור סדרת וקטורי הזיכרון m_t.", 2) | Out-Null

# 6. Paragraph (body 2)
$d.Content.Find.Execute("כמובן שאנו מאבדים כאן מהעוצמה של המנגנון תשומת הלב הרגיל אבל זה יעזור לנו לפתור את סוגיית הסיבוכיות הריבועית בזמן ההיסק. למעשה המחברים מוכיחים (ראו את התמונה למעלה) כי ניתן לממש את המנגנון הזה לסדרתי בעל סיבוכיות לינארית במונחי אורך הקלט. כמובן בזמן האימון ניתן לחשב חיזוי של כמה טוקנים בו זמנית (לפי היכולת החישובית שעומדת לרשותנו) וליהנות מהיתרון של מנגנון תשומת הלב הרגיל.", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר המסוקר מנסח מסגרת מתמטית כללית עבור בעיית ייצוג הזיכרון של פונקצית קלט (u(x בתחום . והנה מתחיל הסיבוך: קודם כל פולינומי Legendre הם מקרה פרטי של פונקציות אורתוגונליות במרחב הילברט (יותר נכון מרחב פונקציונלי L של לבג - המקרה הפרטי של הילברט) המצויד בנוסף בפונקציית מידה mu. אוקיי, מה הדבר הזה אומר בעצם? ממש בגדול זה מרחב של פונקציות שהמכפלה הפנימית ביניהן מוגדרת בתור אינטגרל של מכפלתן תחת מידה mu (במקרה הפשוט ביותר מידה mu שווה ל 1 זהותית ואנו מקבלים אינטגרל Riemann רגיל של המכפלה אבל עבור mu מורכבים יותר כמו Riemann-Stieltjes). פונקציות אורתוגונלית במרחב החמוד הזה מוגדרות בתור אלו שהמכפלה הפנימית שלהן שווה ל 0 (תחת מידה mu). פולינומי Legenge הן אורתוגונליים תחת מידה mu השווה ל-  ב-  ואפס בכל מקום אחר. ", 2) | Out-Null

# 7. Last paragraph (body 3) - replace text
$d.Content.Find.Execute("כלומר יש לנו טרנספורמר (מוחלש כמובן) באימון ו- RNN בהיסק. בהמשך נראה ניתן לשפר את הגישה הזו עם SSMs) state-space models).", $true, $false, $false, $false, $false, $true, 1, $false, "אז נניח שיש לנו N פונקציות אורתוגונליות  במרחבנו החמוד. ועכשיו המטרה היא לתאר את הקלט (u(x ב- על ידי . כלומר אנו רוצים לבנות סכום ממושקל (u*(x של  עם מקדמים מסוימים (שימו לב שעבור t-ים שונים מקבלים וקטורי מקדמים שונים וכך שיש לנו כאן פונקציה וקטורית של המקדמים התלויה ב-t).", 2) | Out-Null

# 8. Append new paragraphs after the last paragraph: empty, text, empty, text, ... empty, text
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "כלומר (u*(x צריך לקרב בצורה טובה את הקלט (u(x (כלומר למזער שגיאה ביניהן ב-). והדיוק מחושב בתור אינטגרל של ההפרש הריבועי בין (u*(x ו- (u(x תחת מידה mu (כאמור היא שווה ל-  ב-  עבור כל x ואפס בכל מקום אחר עבור פולינומי Legendre אבל כמובן קיימות עוד אפשרויות). איך נחשב מקדמים הממזערים את ההפרש הזה? לא כזה מסובך: מקדם i שווה למכפלה פנימית (=אינטגרל) בין פונקציה מספר i לפונקצית קלט u תחת אותה מידה mu. "
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "עכשיו איך כל זה קשור למערכות דינמיות לינאריות החמודות שלנו? מתברר כי מערכת דינמית לינארית שתיארנו בסקירה הקודמת עבור וקטור (m(t מתארת את המקדמים של ייצוג הקלט באמצעות N פולינומי Legendre אורתוגונליים תחת מידה mu שהגדרנו לפני. ו-N זה המימד של וקטור הזיכרון (m(t תחת מידה mu הדורשת קרבה אחידה (=זכרון אחיד) בין u* ו- u ב- . "
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "אם נגדיר מידה mu להיות פונקציה (exp(x-t עבור t נתון, מערכת דינמית לינארית אחרת תתאר לנו מקדמים של פולינומי Laguerre (אורתוגונליים תחת mu הזו). שימו לב שמידה זו מגדירה זיכרון הדועך מעריכית כלומר ככל שעובר הזמן מזמן הנוכחי t, הזיכרון הולך ונהיה מעומעם יותר. "
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "בנוסף המאמר מדבר גם על שיטות דיסקרטיזציה של מערכת דינמית זו וגם דן בקשר בינה לבין RNNs."
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "אוקיי, עכשיו סיכום במשפט אחד של המאמר הדי כבד הזה. המחברים בנו מסגרת מתמטית למידול בעיית הזיכרון של פונקציית קלט שישמש אותנו מאחורי הקלעים לבניית מודלי attention כל הדרך לממבה. "

Write-Output "done"